$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39-99 down to 40-100.
$ws.Rows("39:39").Insert()

# Populate the new row 39 with the new data record.
$ws.Range("A39").Value = 6
$ws.Range("B39").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C39").Value = "Metropolitana"
$ws.Range("D39").Value = 44469
$ws.Range("E39").Value = 13
$ws.Range("F39").Value = 100112001
$ws.Range("G39").Value = "Berenjena"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 190
$ws.Range("K39").Value = 8000
$ws.Range("L39").Value = 9000
$ws.Range("M39").Value = 8579
$ws.Range("N39").Value = "$/caja 60 unidades"
$ws.Range("O39").Value = "Provincia de Huasco"
$ws.Range("P39").Value = 143
$ws.Range("Q39").Value = 60
$ws.Range("R39").Value = "Hortaliza"
